$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.70"
$ws.Range("E2").Value = "'2.08%"
$ws.Range("D3").Value = "'44.07"
$ws.Range("E3").Value = "'6.40%"
$ws.Range("D4").Value = "'5.085"
$ws.Range("E4").Value = "'0.93%"
$ws.Range("D5").Value = "'0.07690"
$ws.Range("E5").Value = "'2.98%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.616"
$ws.Range("E6").Value = "'2.36%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'1.045"
$ws.Range("E7").Value = "'12.64%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1289"
$ws.Range("E8").Value = "'9.02%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1863"
$ws.Range("E9").Value = "'1.47%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.09279"
$ws.Range("E10").Value = "'4.40%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04195"
$ws.Range("E11").Value = "'0.43%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.1047"
$ws.Range("E12").Value = "'-0.35%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001280"
$ws.Range("E13").Value = "'0.29%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.005759"
$ws.Range("E14").Value = "'-1.58%"
$ws.Range("B15").Value = "UpBots"
$ws.Range("C15").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D15").Value = "'0.007489"
$ws.Range("E15").Value = "'1,911.89%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.345"
$ws.Range("E16").Value = "'0.02%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.421"
$ws.Range("E17").Value = "'1.49%"
$ws.Range("D18").Value = "'2.330"
$ws.Range("E18").Value = "'-3.90%"
$ws.Range("D19").Value = "'0.3350"
$ws.Range("E19").Value = "'1.81%"
$ws.Range("D20").Value = "'8.393"
$ws.Range("E20").Value = "'6.49%"
$ws.Range("D21").Value = "'0.1399"
$ws.Range("E21").Value = "'-0.74%"
$ws.Range("D22").Value = "'0.3176"
$ws.Range("E22").Value = "'7.07%"
$ws.Range("D23").Value = "'0.04182"
$ws.Range("E23").Value = "'3.93%"
$ws.Range("D24").Value = "'0.001285"
$ws.Range("E24").Value = "'1.61%"
$ws.Range("D25").Value = "'0.004417"
$ws.Range("E25").Value = "'14.02%"
$ws.Range("D26").Value = "'0.0001351"
$ws.Range("D38").Value = "'0.02493"
$ws.Range("E38").Value = "'4.12%"
$ws.Range("D39").Value = "'0.05297"
$ws.Range("E39").Value = "'1.78%"
$ws.Range("D40").Value = "'0.005936"
$ws.Range("E40").Value = "'-10.21%"
$ws.Range("D41").Value = "'0.007709"
$ws.Range("E41").Value = "'-0.83%"
$ws.Range("D42").Value = "'0.1349"
$ws.Range("E42").Value = "'2.17%"
$ws.Range("D43").Value = "'0.007348"
$ws.Range("E43").Value = "'-0.39%"
$ws.Range("D44").Value = "'0.007557"
$ws.Range("E44").Value = "'5.25%"
$ws.Range("D45").Value = "'0.3006"
$ws.Range("E45").Value = "'-6.32%"
$ws.Range("D46").Value = "'0.00006662"
$ws.Range("E46").Value = "'7.08%"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("D48").Value = "'0.04261"
$ws.Range("E48").Value = "'-7.54%"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E51").Value = "'0.01%"
